$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date values in column C (numbers are Excel date serials)
$ws.Range("C9").Value = 44117
$ws.Range("C19").Value = 44094
$ws.Range("C20").Value = 44123
$ws.Range("C26").Value = 44111
$ws.Range("C28").Value = 44104

# Update the selected cell shown in the sheet view
$ws.Range("C9").Select()
